# 1) Swap the order of slides 7 ("EMC2 Workflow") and 8 ("Linear Designs")
#    so that "Linear Designs" now precedes "EMC2 Workflow" in the deck.
$p = $ppt.ActivePresentation
$moved = $p.Slides.Item(8)
$moved.MoveTo(7)

# 2) On slide 9, merge the two runs "First, " + "we will do some exercises
#    in the 1-BasicEAMs.R script." into a single run of text, keeping the
#    formatting (dirty="0") of the second run.
$s9 = $p.Slides.Item(9)
$shp = $s9.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# The shape auto-fits to its text (spAutoFit); capture its laid-out height
# beforehand so we can restore it after the text edit below (the text
# itself doesn't change length/wrapping, so the box shouldn't move).
$origHeight = [Math]::Round($shp.Height, 4)

$fullText = $tr.Text
$run1Text = "First, "
$run2Text = "we will do some exercises in the 1-BasicEAMs.R script."
$mergedText = "First, we will do some exercises in the 1-BasicEAMs.R script."

$idx1 = $fullText.IndexOf($run1Text)
$idx2 = $fullText.IndexOf($run2Text)

# Rewrite the second run's characters with the merged text (keeps run 2's
# formatting), then clear out the now-duplicated leading run.
$chars2 = $tr.Characters($idx2 + 1, $run2Text.Length)
$chars2.Text = $mergedText

$chars1 = $tr.Characters($idx1 + 1, $run1Text.Length)
$chars1.Text = ""

$shp.Height = $origHeight
